$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GVA")

$ws.Range("B2").Value = 0.1871289562554429
$ws.Range("C2").Value = 1.495060261846345
$ws.Range("D2").Value = 6.487698093947865
$ws.Range("E2").Value = 2.547096011921786
$ws.Range("F2").Value = 2.599990646951396
$ws.Range("G2").Value = 22
